# Fruta / hortaliza, semanal
# A new weekly record is inserted for "Feria Lagunitas de Puerto Montt" / Frutilla,
# pushing the existing rows 192-195 down to 193-196, and the new row 192 is
# populated with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 192; this shifts rows 192-195 down to 193-196
# and keeps their existing data intact.
$ws.Rows(192).Insert()

# Populate the newly inserted row 192 with the new weekly record.
$ws.Range("A192").Value = 4
$ws.Range("B192").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C192").Value = "Los Lagos"
$ws.Range("D192").Value = 44595
$ws.Range("E192").Value = 10
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100101
$ws.Range("H192").Value = "Berries"
$ws.Range("I192").Value = 100112025
$ws.Range("J192").Value = "Frutilla"
$ws.Range("K192").Value = "Sin especificar"
$ws.Range("L192").Value = "Primera"
$ws.Range("M192").Value = 500
$ws.Range("N192").Value = 9000
$ws.Range("O192").Value = 10000
$ws.Range("P192").Value = 9500
$ws.Range("Q192").Value = "`$/caja 7 kilos"
$ws.Range("R192").Value = "Región de La Araucanía"
$ws.Range("S192").Value = 1357
$ws.Range("T192").Value = 7
